# Refresh the crypto price/volume snapshot for Sheet1 (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.075.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.37%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.790.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.40%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'228.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.29%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -1.07%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'32.84"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.49%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.289"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.36%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0714"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.65%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.58%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.049.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.45%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'11.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.74%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.788.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.60%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.625"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.61%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'34.055.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.52%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -3.87%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -1.70%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'245.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.90%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0788"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.62%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.11%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'10.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.42%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.11%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -3.07%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'160.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.11%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'16.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.02%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.59%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.66%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "'  +1.58%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -3.37%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.37%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -2.26%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -3.46%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.402.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.04%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.658"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.49%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.50%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -0.92%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.19%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.43%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.918"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.71%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'78.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.53%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -2.29%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +9.76%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +10.85%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +3.53%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'109.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.23%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.27%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -2.89%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.948.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.08%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.25%  "
$ws.Range("E51").Style = "Normal"
